$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newly-scraped price row (row 24) under the existing tracked
# history (Date, Price, Discount, Incredible). Force text formatting first
# so values like the date and the big price number are stored as literal
# strings (shared strings), matching how the rest of the sheet is written,
# instead of being auto-converted to a date serial / number by Excel.
$newRow = $ws.Range("A24:D24")
$newRow.NumberFormat = "@"

$ws.Range("A24").Value = "2026-02-07"
$ws.Range("B24").Value = "1030000"
$ws.Range("C24").Value = "0"
$ws.Range("D24").Value = "0"

# Drop the temporary text-format override so the new cells end up with the
# same default styling as every other row in the sheet.
$newRow.Style = "Normal"
